# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal"):
# insert two new daily records (Lechuga "Conconina(o)" and "Escarola",
# Primera quality) for the Terminal Hortofrutícola Agro Chillán market on
# 2023-04-05 (serial date 45021), pushing all later rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 985.
$ws.Rows("985:986").Insert()

# --- New row 985: Lechuga / Conconina(o) / Primera ---
$ws.Cells.Item(985, 1).Value2  = 7
$ws.Cells.Item(985, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(985, 3).Value2  = "Ñuble"
$ws.Cells.Item(985, 4).Value2  = 45021
$ws.Cells.Item(985, 5).Value2  = 16
$ws.Cells.Item(985, 6).Value2  = 100112033
$ws.Cells.Item(985, 7).Value2  = "Lechuga"
$ws.Cells.Item(985, 8).Value2  = "Conconina(o)"
$ws.Cells.Item(985, 9).Value2  = "Primera"
$ws.Cells.Item(985, 10).Value2 = 170
$ws.Cells.Item(985, 11).Value2 = 6000
$ws.Cells.Item(985, 12).Value2 = 7000
$ws.Cells.Item(985, 13).Value2 = 6588
$ws.Cells.Item(985, 14).Value2 = "`$/caja 10 unidades"
$ws.Cells.Item(985, 15).Value2 = "Región del Maule"
$ws.Cells.Item(985, 16).Value2 = 659
$ws.Cells.Item(985, 17).Value2 = 10
$ws.Cells.Item(985, 18).Value2 = "Hortaliza"

# --- New row 986: Lechuga / Escarola / Primera ---
$ws.Cells.Item(986, 1).Value2  = 7
$ws.Cells.Item(986, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(986, 3).Value2  = "Ñuble"
$ws.Cells.Item(986, 4).Value2  = 45021
$ws.Cells.Item(986, 5).Value2  = 16
$ws.Cells.Item(986, 6).Value2  = 100112033
$ws.Cells.Item(986, 7).Value2  = "Lechuga"
$ws.Cells.Item(986, 8).Value2  = "Escarola"
$ws.Cells.Item(986, 9).Value2  = "Primera"
$ws.Cells.Item(986, 10).Value2 = 220
$ws.Cells.Item(986, 11).Value2 = 8000
$ws.Cells.Item(986, 12).Value2 = 9000
$ws.Cells.Item(986, 13).Value2 = 8318
$ws.Cells.Item(986, 14).Value2 = "`$/caja 15 unidades"
$ws.Cells.Item(986, 15).Value2 = "Región del Maule"
$ws.Cells.Item(986, 16).Value2 = 555
$ws.Cells.Item(986, 17).Value2 = 15
$ws.Cells.Item(986, 18).Value2 = "Hortaliza"
